# Applies the "first code review" changes to the Settings sheet:
#  - inserts a new "Carturesti_Code" setting row above the
#    CarturestiDT_Columns row (old row 14)
#  - inserts a new "eMAG_Code" setting row above the eMAG_Sheet row
#    (old row 24, new row 26 after the first insert)
#  - fixes a typo in the eMAGDT_Columns value (double space -> single)
#  - bumps DelayForFilters from 1 to 1.5
#  - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Fix the eMAGDT_Columns value typo (extra space removed), before the
#     row inserts below shift it out from under row 25 ---
$ws.Range("B25").Value = "Title,Author,Overall review,No Reviews,Price"

# --- Insert the Carturesti_Code row block (2 rows) above old row 14 ---
$ws.Rows("14:15").Insert()
$ws.Rows("14:15").RowHeight = 14.25

$ws.Range("A14").Value = "Carturesti_Code"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = "This will be used as reference when uploading queue items."

# --- Insert the eMAG_Code row above the eMAG_Sheet row (now row 26) ---
$ws.Rows("26:26").Insert()
$ws.Rows("26:26").RowHeight = 14.25

$ws.Range("A26").Value = "eMAG_Code"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "This will be used as reference whenn uploading queue items."

# --- DelayForFilters value bump ---
$ws.Range("B34").Value = 1.5

# --- Update selection to match the saved view state ---
$ws.Range("C26").Select()
